$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '41.767.22'
Set-TextValue "E2" '  +0.59%  '
Set-TextValue "D3" '2.478.53'
Set-TextValue "E3" '  +0.43%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.04%  '
Set-TextValue "D5" '319.44'
Set-TextValue "E5" '  +1.53%  '
Set-TextValue "D6" '93.46'
Set-TextValue "E6" '  +1.78%  '
Set-TextValue "E7" '  +0.88%  '
Set-TextValue "E8" '  -0.02%  '
Set-TextValue "E9" '  +1.77%  '
Set-TextValue "D10" '0.0870'
Set-TextValue "E10" '  +10.04%  '
Set-TextValue "D11" '33.34'
Set-TextValue "E11" '  +3.57%  '
Set-TextValue "E12" '  +0.70%  '
Set-TextValue "D13" '2.861.60'
Set-TextValue "E13" '  +0.50%  '
Set-TextValue "E14" '  +1.26%  '
Set-TextValue "D15" '15.87'
Set-TextValue "E15" '  -0.42%  '
Set-TextValue "D16" '2.489.53'
Set-TextValue "E16" '  +0.61%  '
Set-TextValue "E17" '  +2.95%  '
Set-TextValue "D18" '41.728.87'
Set-TextValue "E18" '  +0.54%  '
Set-TextValue "D19" '6.49'
Set-TextValue "E19" '  +0.43%  '
Set-TextValue "D20" '0.0₃0953'
Set-TextValue "E20" '  +1.20%  '
Set-TextValue "D21" '71.33'
Set-TextValue "E21" '  +0.18%  '
Set-TextValue "D22" '11.35'
Set-TextValue "E22" '  +2.68%  '
Set-TextValue "D23" '240.12'
Set-TextValue "E23" '  +1.94%  '
Set-TextValue "E24" '  +1.54%  '
Set-TextValue "E25" '  +2.71%  '
Set-TextValue "E26" '  +0.05%  '
Set-TextValue "D27" '24.79'
Set-TextValue "E27" '  +0.71%  '
Set-TextValue "E28" '  +1.51%  '
Set-TextValue "D29" '9.85'
Set-TextValue "E29" '  +1.91%  '
Set-TextValue "D30" '36.27'
Set-TextValue "E30" '  +2.58%  '
Set-TextValue "D31" '158.33'
Set-TextValue "E31" '  +1.71%  '
Set-TextValue "E32" '  +2.03%  '
Set-TextValue "E33" '  -0.14%  '
Set-TextValue "B34" 'WEMIXToken'
Set-TextValue "C34" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D34" '2.59'
Set-TextValue "E34" '  +1.13%  '
Set-TextValue "B35" 'Hedera'
Set-TextValue "C35" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D35" '0.0768'
Set-TextValue "E35" '  +1.49%  '
Set-TextValue "D36" '17.44'
Set-TextValue "E36" '  +0.77%  '
Set-TextValue "D37" '1.88'
Set-TextValue "E37" '  +6.20%  '
Set-TextValue "E38" '  +2.80%  '
Set-TextValue "E39" '  +1.74%  '
Set-TextValue "E40" '  +0.40%  '
Set-TextValue "D41" '4.04'
Set-TextValue "E41" '  +0.40%  '
Set-TextValue "E42" '  +10.47%  '
Set-TextValue "D43" '1.990.84'
Set-TextValue "E43" '  +2.60%  '
Set-TextValue "D44" '19.12'
Set-TextValue "E44" '  +4.32%  '
Set-TextValue "E45" '  +1.30%  '
Set-TextValue "D46" '3.00'
Set-TextValue "E46" '  +2.67%  '
Set-TextValue "D47" '9.45'
Set-TextValue "E47" '  +5.11%  '
Set-TextValue "D48" '2.718.45'
Set-TextValue "E48" '  +0.40%  '
Set-TextValue "D49" '97.75'
Set-TextValue "E49" '  +1.12%  '
Set-TextValue "D50" '74.44'
Set-TextValue "E50" '  +3.60%  '
Set-TextValue "D51" '67.35'
Set-TextValue "E51" '  +0.64%  '
